# Add results and unfolding with 100 keV threshold
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Beta" row (row 2) values ---
$ws.Range("C2").Value = 40.35273080373261
$ws.Range("E2").Value = 0.05768901818751787
$ws.Range("F2").Value = 249.4604156738437
$ws.Range("G2").Value = 1.452772998844936
$ws.Range("H2").Value = 520.9915530125774
$ws.Range("I2").Value = 0.004182014977651938
$ws.Range("J2").Value = 0.000005483887722717556
$ws.Range("K2").Value = 0.01208460435650103
$ws.Range("L2").Value = 0.0866379204398133
$ws.Range("M2").Value = 0.0002766858027195116
$ws.Range("N2").Value = 0.2002514646196779

# --- Update existing "Gamma" row (row 3) values ---
$ws.Range("F3").Value = 0.0001726308312382031
$ws.Range("G3").Value = 0.00007030986138709599
$ws.Range("H3").Value = 0.000285513147140896
$ws.Range("I3").Value = 0.0001601665874816452
$ws.Range("J3").Value = 0.00006567929002917165
$ws.Range("K3").Value = 0.0002635850885777029
$ws.Range("L3").Value = 0.0001788441461216513
$ws.Range("M3").Value = 0.00007272893004767802
$ws.Range("N3").Value = 0.0002958506245513513

# --- Add new "Beta + Gamma" row (row 4) ---
# Copy formatting from row 3 (A3 has the bordered/bold/centered style)
$ws.Range("A3:N3").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.35273080373261
$ws.Range("D4").Value = 0.0007432820064133916
$ws.Range("E4").Value = 0.05768901818751787
$ws.Range("F4").Value = 249.460588304675
$ws.Range("G4").Value = 1.452843308706323
$ws.Range("H4").Value = 520.9918385257246
$ws.Range("I4").Value = 0.004342181565133585
$ws.Range("J4").Value = 0.00007116317775188921
$ws.Range("K4").Value = 0.01234818944507873
$ws.Range("L4").Value = 0.08681676458593496
$ws.Range("M4").Value = 0.0003494147327671897
$ws.Range("N4").Value = 0.2005473152442293
